$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023505114265091
$ws.Cells.Item(2, 4).Value = 1.045302441635666
$ws.Cells.Item(2, 5).Value = 1.02411664715054
$ws.Cells.Item(2, 6).Value = 1.049197042350438
$ws.Cells.Item(2, 9).Value = 1.038398395490394
$ws.Cells.Item(2, 10).Value = 1.02868495775162
$ws.Cells.Item(2, 11).Value = 1.0480710626297
$ws.Cells.Item(2, 12).Value = 1.026946091997175
$ws.Cells.Item(2, 13).Value = 1.051954771982741
$ws.Cells.Item(2, 14).Value = 1.013620620313449

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.024289244602734
$ws.Cells.Item(3, 4).Value = 1.045910744445182
$ws.Cells.Item(3, 5).Value = 1.024777724456729
$ws.Cells.Item(3, 6).Value = 1.049983800740937
$ws.Cells.Item(3, 9).Value = 1.038582137413128
$ws.Cells.Item(3, 10).Value = 1.02910845225797
$ws.Cells.Item(3, 11).Value = 1.048491307656755
$ws.Cells.Item(3, 12).Value = 1.02741458972594
$ws.Cells.Item(3, 13).Value = 1.052553793851092
$ws.Cells.Item(3, 14).Value = 1.013761410948608

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02479725986983
$ws.Cells.Item(4, 4).Value = 1.046304641019365
$ws.Cells.Item(4, 5).Value = 1.025206427292934
$ws.Cells.Item(4, 6).Value = 1.050493708345829
$ws.Cells.Item(4, 9).Value = 1.038699885011438
$ws.Cells.Item(4, 10).Value = 1.02938244571725
$ws.Cells.Item(4, 11).Value = 1.048762822640324
$ws.Cells.Item(4, 12).Value = 1.027717997817982
$ws.Cells.Item(4, 13).Value = 1.052941564170229
$ws.Cells.Item(4, 14).Value = 1.013852475300868

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025010978640516
$ws.Cells.Item(5, 4).Value = 1.046470301110373
$ws.Cells.Item(5, 5).Value = 1.025386877482501
$ws.Cells.Item(5, 6).Value = 1.050708268174184
$ws.Cells.Item(5, 9).Value = 1.038749110904229
$ws.Cells.Item(5, 10).Value = 1.02949762271141
$ws.Cells.Item(5, 11).Value = 1.048876867470713
$ws.Cells.Item(5, 12).Value = 1.027845611109209
$ws.Cells.Item(5, 13).Value = 1.05310462002037
$ws.Cells.Item(5, 14).Value = 1.013890749532217

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025046871656822
$ws.Cells.Item(6, 4).Value = 1.046498119978525
$ws.Cells.Item(6, 5).Value = 1.025417188914988
$ws.Cells.Item(6, 6).Value = 1.05074430506527
$ws.Cells.Item(6, 9).Value = 1.038757359984012
$ws.Cells.Item(6, 10).Value = 1.029516960820199
$ws.Cells.Item(6, 11).Value = 1.048896010187274
$ws.Cells.Item(6, 12).Value = 1.027867041450706
$ws.Cells.Item(6, 13).Value = 1.053131999943377
$ws.Cells.Item(6, 14).Value = 1.0138971753899

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024800115006681
$ws.Cells.Item(7, 4).Value = 1.046306854320927
$ws.Cells.Item(7, 5).Value = 1.025208837601332
$ws.Cells.Item(7, 6).Value = 1.050496574542592
$ws.Cells.Item(7, 9).Value = 1.038700543852706
$ws.Cells.Item(7, 10).Value = 1.029383984758224
$ws.Cells.Item(7, 11).Value = 1.048764346908041
$ws.Cells.Item(7, 12).Value = 1.027719702756962
$ws.Cells.Item(7, 13).Value = 1.052943742786594
$ws.Cells.Item(7, 14).Value = 1.013852986759836

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023769983155966
$ws.Cells.Item(8, 4).Value = 1.045507960408401
$ws.Cells.Item(8, 5).Value = 1.024339865208232
$ws.Cells.Item(8, 6).Value = 1.049462759448068
$ws.Cells.Item(8, 9).Value = 1.038460728612568
$ws.Cells.Item(8, 10).Value = 1.028828086344501
$ws.Cells.Item(8, 11).Value = 1.048213170726473
$ws.Cells.Item(8, 12).Value = 1.02710436857707
$ws.Cells.Item(8, 13).Value = 1.052157179346933
$ws.Cells.Item(8, 14).Value = 1.013668208470619

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021959667431124
$ws.Cells.Item(9, 4).Value = 1.044102471162545
$ws.Cells.Item(9, 5).Value = 1.0228159207287
$ws.Cells.Item(9, 6).Value = 1.047647440057659
$ws.Cells.Item(9, 9).Value = 1.038029409045867
$ws.Cells.Item(9, 10).Value = 1.027848304861281
$ws.Cells.Item(9, 11).Value = 1.047238842221121
$ws.Cells.Item(9, 12).Value = 1.026022115124723
$ws.Cells.Item(9, 13).Value = 1.05077248473013
$ws.Cells.Item(9, 14).Value = 1.013342345479961

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02075619721681
$ws.Cells.Item(10, 4).Value = 1.043167120626321
$ws.Cells.Item(10, 5).Value = 1.021804977003331
$ws.Cells.Item(10, 6).Value = 1.046441657623225
$ws.Cells.Item(10, 9).Value = 1.037736044489974
$ws.Cells.Item(10, 10).Value = 1.027195047625016
$ws.Cells.Item(10, 11).Value = 1.046587305857369
$ws.Cells.Item(10, 12).Value = 1.025302071371045
$ws.Cells.Item(10, 13).Value = 1.049850359566175
$ws.Cells.Item(10, 14).Value = 1.013124956056312

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020235911968916
$ws.Cells.Item(11, 4).Value = 1.042762518487121
$ws.Cells.Item(11, 5).Value = 1.021368441849259
$ws.Cells.Item(11, 6).Value = 1.045920619078857
$ws.Cells.Item(11, 9).Value = 1.03760764587245
$ws.Cells.Item(11, 10).Value = 1.026912179715096
$ws.Cells.Item(11, 11).Value = 1.046304731982087
$ws.Cells.Item(11, 12).Value = 1.0249906482566
$ws.Cells.Item(11, 13).Value = 1.049451330709731
$ws.Cells.Item(11, 14).Value = 1.013030794859287

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020042780674012
$ws.Cells.Item(12, 4).Value = 1.042612295022892
$ws.Cells.Item(12, 5).Value = 1.021206477036328
$ws.Cells.Item(12, 6).Value = 1.045727245652667
$ws.Cells.Item(12, 9).Value = 1.037559748017292
$ws.Cells.Item(12, 10).Value = 1.026807110801699
$ws.Cells.Item(12, 11).Value = 1.046199704748214
$ws.Cells.Item(12, 12).Value = 1.02487502746485
$ws.Cells.Item(12, 13).Value = 1.049303154182309
$ws.Cells.Item(12, 14).Value = 1.012995815137228

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020084202269875
$ws.Cells.Item(13, 4).Value = 1.042644515546398
$ws.Cells.Item(13, 5).Value = 1.021241210690506
$ws.Cells.Item(13, 6).Value = 1.045768717474697
$ws.Cells.Item(13, 9).Value = 1.037570031526989
$ws.Cells.Item(13, 10).Value = 1.026829648382085
$ws.Cells.Item(13, 11).Value = 1.046222236440648
$ws.Cells.Item(13, 12).Value = 1.024899825977097
$ws.Cells.Item(13, 13).Value = 1.049334936684984
$ws.Cells.Item(13, 14).Value = 1.013003318582619

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020219945098775
$ws.Cells.Item(14, 4).Value = 1.042750099655585
$ws.Cells.Item(14, 5).Value = 1.021355050026871
$ws.Cells.Item(14, 6).Value = 1.045904631423003
$ws.Cells.Item(14, 9).Value = 1.037603690795267
$ws.Cells.Item(14, 10).Value = 1.026903494661599
$ws.Cells.Item(14, 11).Value = 1.046296051750223
$ws.Cells.Item(14, 12).Value = 1.024981089862921
$ws.Cells.Item(14, 13).Value = 1.049439081553537
$ws.Cells.Item(14, 14).Value = 1.013027903504691

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.020303597461508
$ws.Cells.Item(15, 4).Value = 1.042815162023801
$ws.Cells.Item(15, 5).Value = 1.021425214610628
$ws.Cells.Item(15, 6).Value = 1.045988394223367
$ws.Cells.Item(15, 9).Value = 1.03762440225843
$ws.Cells.Item(15, 10).Value = 1.026948993939752
$ws.Cells.Item(15, 11).Value = 1.046341523004392
$ws.Cells.Item(15, 12).Value = 1.025031166638294
$ws.Cells.Item(15, 13).Value = 1.049503254074643
$ws.Cells.Item(15, 14).Value = 1.01304305056244

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020790744378598
$ws.Cells.Item(16, 4).Value = 1.043193981569543
$ws.Cells.Item(16, 5).Value = 1.021833974068419
$ws.Cells.Item(16, 6).Value = 1.046476260038885
$ws.Cells.Item(16, 9).Value = 1.037744537100268
$ws.Cells.Item(16, 10).Value = 1.027213820680364
$ws.Cells.Item(16, 11).Value = 1.046606049914384
$ws.Cells.Item(16, 12).Value = 1.02532274721801
$ws.Cells.Item(16, 13).Value = 1.049876847381511
$ws.Cells.Item(16, 14).Value = 1.013131204629714

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021096541009372
$ws.Cells.Item(17, 4).Value = 1.043431716508365
$ws.Cells.Item(17, 5).Value = 1.022090703459552
$ws.Cells.Item(17, 6).Value = 1.046782574293826
$ws.Cells.Item(17, 9).Value = 1.037819528412602
$ws.Cells.Item(17, 10).Value = 1.027379939559618
$ws.Cells.Item(17, 11).Value = 1.046771860109198
$ws.Cells.Item(17, 12).Value = 1.025505745571573
$ws.Cells.Item(17, 13).Value = 1.050111262841604
$ws.Cells.Item(17, 14).Value = 1.013186493607516

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.021274986517435
$ws.Cells.Item(18, 4).Value = 1.043570422775539
$ws.Cells.Item(18, 5).Value = 1.022240566000028
$ws.Cells.Item(18, 6).Value = 1.046961345614385
$ws.Cells.Item(18, 9).Value = 1.037863137346693
$ws.Cells.Item(18, 10).Value = 1.027476833493578
$ws.Cells.Item(18, 11).Value = 1.046868530450914
$ws.Cells.Item(18, 12).Value = 1.025612520109271
$ws.Cells.Item(18, 13).Value = 1.050248018125191
$ws.Cells.Item(18, 14).Value = 1.013218739793267

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.021335845254645
$ws.Cells.Item(19, 4).Value = 1.043617724665067
$ws.Cells.Item(19, 5).Value = 1.022291684976081
$ws.Cells.Item(19, 6).Value = 1.047022319484402
$ws.Cells.Item(19, 9).Value = 1.03787798441015
$ws.Cells.Item(19, 10).Value = 1.027509871713863
$ws.Cells.Item(19, 11).Value = 1.046901485039351
$ws.Cells.Item(19, 12).Value = 1.025648933308042
$ws.Cells.Item(19, 13).Value = 1.050294652277077
$ws.Cells.Item(19, 14).Value = 1.013229734396454

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021063723688008
$ws.Cells.Item(20, 4).Value = 1.043406205698813
$ws.Cells.Item(20, 5).Value = 1.022063146747744
$ws.Cells.Item(20, 6).Value = 1.046749698962011
$ws.Cells.Item(20, 9).Value = 1.03781149622507
$ws.Cells.Item(20, 10).Value = 1.027362116607731
$ws.Cells.Item(20, 11).Value = 1.046754074778255
$ws.Cells.Item(20, 12).Value = 1.025486107981107
$ws.Cells.Item(20, 13).Value = 1.050086109711599
$ws.Cells.Item(20, 14).Value = 1.013180561924129

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020179968731315
$ws.Cells.Item(21, 4).Value = 1.042719005989702
$ws.Cells.Item(21, 5).Value = 1.021321522090664
$ws.Cells.Item(21, 6).Value = 1.045864603621648
$ws.Cells.Item(21, 9).Value = 1.0375937846284
$ws.Cells.Item(21, 10).Value = 1.026881748739649
$ws.Cells.Item(21, 11).Value = 1.046274316810632
$ws.Cells.Item(21, 12).Value = 1.024957158127038
$ws.Cells.Item(21, 13).Value = 1.049408412350832
$ws.Cells.Item(21, 14).Value = 1.013020663966254

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.019625045724756
$ws.Cells.Item(22, 4).Value = 1.042287306317671
$ws.Cells.Item(22, 5).Value = 1.020856297499817
$ws.Cells.Item(22, 6).Value = 1.045309056057573
$ws.Cells.Item(22, 9).Value = 1.03745571571713
$ws.Cells.Item(22, 10).Value = 1.026579727802773
$ws.Cells.Item(22, 11).Value = 1.045972288704533
$ws.Cells.Item(22, 12).Value = 1.0246249088686
$ws.Cells.Item(22, 13).Value = 1.048982552699933
$ws.Cells.Item(22, 14).Value = 1.012920106509282

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.019919151218925
$ws.Cells.Item(23, 4).Value = 1.042516122716226
$ws.Cells.Item(23, 5).Value = 1.021102820321341
$ws.Cells.Item(23, 6).Value = 1.045603471816485
$ws.Cells.Item(23, 9).Value = 1.037529020678843
$ws.Cells.Item(23, 10).Value = 1.026739833867433
$ws.Cells.Item(23, 11).Value = 1.04613243551892
$ws.Cells.Item(23, 12).Value = 1.024801009434615
$ws.Cells.Item(23, 13).Value = 1.049208286000809
$ws.Cells.Item(23, 14).Value = 1.01297341596911

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021078552177496
$ws.Cells.Item(24, 4).Value = 1.043417732813776
$ws.Cells.Item(24, 5).Value = 1.022075598079542
$ws.Cells.Item(24, 6).Value = 1.04676455359017
$ws.Cells.Item(24, 9).Value = 1.037815126033556
$ws.Cells.Item(24, 10).Value = 1.027370170033907
$ws.Cells.Item(24, 11).Value = 1.046762111339713
$ws.Cells.Item(24, 12).Value = 1.02549498125594
$ws.Cells.Item(24, 13).Value = 1.050097475252234
$ws.Cells.Item(24, 14).Value = 1.013183242205669

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022427084777717
$ws.Cells.Item(25, 4).Value = 1.04446554313846
$ws.Cells.Item(25, 5).Value = 1.023209020552807
$ws.Cells.Item(25, 6).Value = 1.048115972303808
$ws.Cells.Item(25, 9).Value = 1.03814194496918
$ws.Cells.Item(25, 10).Value = 1.02810161972978
$ws.Cells.Item(25, 11).Value = 1.04749108581684
$ws.Cells.Item(25, 12).Value = 1.02630165266805
$ws.Cells.Item(25, 13).Value = 1.051130292454154
$ws.Cells.Item(25, 14).Value = 1.013426616863415

Write-Host "done"
